$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jarno")

$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 44988
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Palvelutasot"
